$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.456.97'
$ws.Range("E2").Value = '  -2.00%  '
$ws.Range("D3").Value = '3.503.40'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.05'
$ws.Range("E5").Value = '  -2.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.12'
$ws.Range("E6").Value = '  -3.24%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.506.31'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("E9").Value = '  -3.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.133'
$ws.Range("E10").Value = '  -5.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.88'
$ws.Range("E11").Value = '  -2.03%  '
$ws.Range("E12").Value = '  -3.91%  '
$ws.Range("D13").Value = '4.120.10'
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '30.33'
$ws.Range("E14").Value = '  -6.07%  '
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = '66.511.83'
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000172'
$ws.Range("E17").Value = '  -3.93%  '
$ws.Range("D18").Value = '3.518.08'
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("E19").Value = '  -5.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.84'
$ws.Range("E20").Value = '  -3.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '379.00'
$ws.Range("E21").Value = '  -3.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.85'
$ws.Range("E22").Value = '  -1.85%  '
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.75'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '72.13'
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000120'
$ws.Range("E27").Value = '  -2.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.86'
$ws.Range("E28").Value = '  -5.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.174'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '24.44'
$ws.Range("E31").Value = '  +3.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.89'
$ws.Range("E32").Value = '  -4.75%  '
$ws.Range("E33").Value = '  -3.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.32'
$ws.Range("E34").Value = '  -7.37%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -3.03%  '
$ws.Range("E37").Value = '  -2.85%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.38'
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '29.38'
$ws.Range("E39").Value = '  +10.76%  '
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.78'
$ws.Range("E41").Value = '  -6.30%  '
$ws.Range("E42").Value = '  -2.94%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.45'
$ws.Range("E43").Value = '  -5.88%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.55'
$ws.Range("E44").Value = '  -10.67%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0696'
$ws.Range("E45").Value = '  -4.61%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.681.96'
$ws.Range("E46").Value = '  -5.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.71'
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.58'
$ws.Range("E48").Value = '  -8.88%  '
$ws.Range("E49").Value = '  -3.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '313.50'
$ws.Range("E50").Value = '  -7.85%  '
$ws.Range("E51").Value = '  -5.30%  '
